# Update the header date line.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-09-04 Thursday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-09-05 Friday", 2)

# Update the division-problem answers laid out in the single 5-column table.
# Only every 4th row (1, 5, 9, 13, 17) actually holds text; the others are blank
# spacer rows. Each cell's text is replaced in place by (row, column) position
# so there is no risk of accidental re-matching between old/new values that
# happen to coincide (e.g. "49÷6=8, 1" is both an old value in row 1 and a new
# value in row 17).
$t = $d.Tables.Item(1)

$rowsData = @(
    @(1,  @("75÷3=25, 0", "13÷3=4, 1", "40÷9=4, 4", "92÷4=23, 0", "11÷8=1, 3")),
    @(5,  @("75÷7=10, 5", "21÷7=3, 0", "42÷5=8, 2", "45÷3=15, 0", "38÷5=7, 3")),
    @(9,  @("89÷2=44, 1", "64÷2=32, 0", "94÷8=11, 6", "97÷7=13, 6", "50÷5=10, 0")),
    @(13, @("77÷3=25, 2", "59÷2=29, 1", "15÷2=7, 1", "37÷4=9, 1", "18÷5=3, 3")),
    @(17, @("17÷5=3, 2", "14÷8=1, 6", "49÷6=8, 1", "26÷5=5, 1", "43÷4=10, 3"))
)

foreach ($rowEntry in $rowsData) {
    $rowIndex = $rowEntry[0]
    $values = $rowEntry[1]
    for ($col = 1; $col -le $values.Count; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
